# Generate Report for Handoff
# Adds a new row (row 3) to each of the three worksheets (Overview, zh-cn,
# de-de) describing the newly-handed-off file
# "23a569ff-5c93-42a2-aee3-6bb02e9f9d3a...md", mirroring the existing row 2
# which describes "9d77cfaa-2b44-466c-a341-0e8308d097c2...md".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Shared literal values (kept identical across sheets)
# ---------------------------------------------------------------------
$newMdFile = '23a569ff-5c93-42a2-aee3-6bb02e9f9d3aooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$newMdDisplay = 'e2e\23a569ff-5c93-42a2-aee3-6bb02e9f9d3aooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$newZhXlf = '23a569ff-5c93-42a2-aee3-6bb02e9f9d3aoooooooooooooooooooooooooooooooooooooooo.d6f4eb904ab7916d515852f8abe20c6dd0854e6e.zh-cn.xlf'
$newDeXlf = '23a569ff-5c93-42a2-aee3-6bb02e9f9d3aoooooooooooooooooooooooooooooooooooooooo.d6f4eb904ab7916d515852f8abe20c6dd0854e6e.de-de.xlf'
$hoCommit = 'ae660cccacb387e63e997b29901a839f9c78b097'
$mdUrl = 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/' + $hoCommit + '/e2e/' + $newMdFile

$statusText = 'Ready for handoff'
$zhHandoffTime = '2016-09-06 16:53:19'
$deHandoffTime = '2016-09-06 16:53:24'
$dateFmt = 'yyyy-mm-dd HH:mm:ss'

# ---------------------------------------------------------------------
# Overview sheet: append summary row 3
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item('Overview')
$loOv = $wsOv.ListObjects.Item(1)
$loOv.ListRows.Add() | Out-Null

$wsOv.Range('A3').Value = $newMdFile
$wsOv.Range('B3').Value = $newMdDisplay
$wsOv.Range('C3').Value = '.md'
$wsOv.Range('E3').Value = $statusText
$wsOv.Range('F3').Value = $statusText
$wsOv.Range('G3').Value = $deHandoffTime
$wsOv.Range('G3').NumberFormat = $dateFmt

$wsOv.Hyperlinks.Add($wsOv.Range('B3'), $mdUrl, '', '', $newMdDisplay) | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet: append detail row 3
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item('zh-cn')
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range('A3').Value = $newMdFile
$wsZh.Range('B3').Value = '.md'
$wsZh.Range('C3').Value = $statusText
$wsZh.Range('D3').Value = 'e2e'
$wsZh.Range('E3').Value = 'ht'
$wsZh.Range('F3').Value = "'False"
$wsZh.Range('G3').Value = $newZhXlf
$wsZh.Range('H3').Value = $zhHandoffTime
$wsZh.Range('H3').NumberFormat = $dateFmt
$wsZh.Range('K3').Value = '0001-01-01 00:00:00'
$wsZh.Range('K3').NumberFormat = $dateFmt
$wsZh.Range('M3').Value = "'True"
$wsZh.Range('O3').Value = "'False"

$wsZh.Hyperlinks.Add($wsZh.Range('A3'), $mdUrl, '', '', $newMdDisplay) | Out-Null

# ---------------------------------------------------------------------
# de-de sheet: append detail row 3
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item('de-de')
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range('A3').Value = $newMdFile
$wsDe.Range('B3').Value = '.md'
$wsDe.Range('C3').Value = $statusText
$wsDe.Range('D3').Value = 'e2e'
$wsDe.Range('E3').Value = 'ht'
$wsDe.Range('F3').Value = "'False"
$wsDe.Range('G3').Value = $newDeXlf
$wsDe.Range('H3').Value = $deHandoffTime
$wsDe.Range('H3').NumberFormat = $dateFmt
$wsDe.Range('K3').Value = '0001-01-01 00:00:00'
$wsDe.Range('K3').NumberFormat = $dateFmt
$wsDe.Range('M3').Value = "'True"
$wsDe.Range('O3').Value = "'False"

$wsDe.Hyperlinks.Add($wsDe.Range('A3'), $mdUrl, '', '', $newMdDisplay) | Out-Null

Write-Host 'Generate Report for Handoff: done'
